$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 470, shifting existing rows 470-490 down to 471-491
$ws.Rows.Item(470).Insert()

# Populate the newly inserted row 470 with the new record
$ws.Cells.Item(470, 1).Value = 5
$ws.Cells.Item(470, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(470, 3).Value = "Maule"
$ws.Cells.Item(470, 4).Value = 45147
$ws.Cells.Item(470, 4).NumberFormat = $ws.Cells.Item(471, 4).NumberFormat
$ws.Cells.Item(470, 5).Value = 7
$ws.Cells.Item(470, 6).Value = "Fruta"
$ws.Cells.Item(470, 7).Value = 100101
$ws.Cells.Item(470, 8).Value = "Berries"
$ws.Cells.Item(470, 9).Value = 100101007
$ws.Cells.Item(470, 10).Value = "Kiwi"
$ws.Cells.Item(470, 11).Value = "Hayward"
$ws.Cells.Item(470, 12).Value = "Primera"
$ws.Cells.Item(470, 13).Value = 200
$ws.Cells.Item(470, 14).Value = 13000
$ws.Cells.Item(470, 15).Value = 13000
$ws.Cells.Item(470, 16).Value = 13000
$ws.Cells.Item(470, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(470, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(470, 19).Value = 722
$ws.Cells.Item(470, 20).Value = 18
